$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 50 (pushes old rows
# 50..108 down to 51..109) and populate it with the new weekly observation.
$ws.Rows.Item(50).Insert()

$ws.Range("A50").Value = 11
$ws.Range("B50").Value = "Vega Monumental Concepción"
$ws.Range("C50").Value = "Bíobío"
$ws.Range("D50").Value = 44790
$ws.Range("E50").Value = 8
$ws.Range("F50").Value = 100112001
$ws.Range("G50").Value = "Berenjena"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 180
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 11000
$ws.Range("M50").Value = 10444
$ws.Range("N50").Value = "$/caja 60 unidades"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 174
$ws.Range("Q50").Value = 60
$ws.Range("R50").Value = "Hortaliza"

Write-Output "done"
